$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 44:45, shifting old rows 44+ down to 46+ (old row 128 -> new row 130).
$ws.Rows("44:45").Insert()

# The old rows 42 and 43 (Murcott data) are still sitting at 42:43 (above the insertion point).
# Copy them down into the freshly-inserted 44:45 slots so they land where the diff expects them.
$ws.Rows("42:43").Copy()
$ws.Rows("44").PasteSpecial()

# Now overwrite rows 42 and 43 with the brand-new weekly data (2022-08-02 / Clemenuless).
$ws.Cells.Item(42, 1).Value = 11
$ws.Cells.Item(42, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(42, 3).Value = "Bíobío"
$ws.Cells.Item(42, 4).Value = [DateTime]"2022-08-02"
$ws.Cells.Item(42, 5).Value = 8
$ws.Cells.Item(42, 6).Value = "Fruta"
$ws.Cells.Item(42, 7).Value = 100102
$ws.Cells.Item(42, 8).Value = "Cítricos"
$ws.Cells.Item(42, 9).Value = 100102004
$ws.Cells.Item(42, 10).Value = "Mandarina"
$ws.Cells.Item(42, 11).Value = "Clemenuless"
$ws.Cells.Item(42, 12).Value = "Primera"
$ws.Cells.Item(42, 13).Value = 100
$ws.Cells.Item(42, 14).Value = 8000
$ws.Cells.Item(42, 15).Value = 9000
$ws.Cells.Item(42, 16).Value = 8500
$ws.Cells.Item(42, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(42, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(42, 19).Value = 472
$ws.Cells.Item(42, 20).Value = 18

$ws.Cells.Item(43, 1).Value = 11
$ws.Cells.Item(43, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(43, 3).Value = "Bíobío"
$ws.Cells.Item(43, 4).Value = [DateTime]"2022-08-02"
$ws.Cells.Item(43, 5).Value = 8
$ws.Cells.Item(43, 6).Value = "Fruta"
$ws.Cells.Item(43, 7).Value = 100102
$ws.Cells.Item(43, 8).Value = "Cítricos"
$ws.Cells.Item(43, 9).Value = 100102004
$ws.Cells.Item(43, 10).Value = "Mandarina"
$ws.Cells.Item(43, 11).Value = "Clemenuless"
$ws.Cells.Item(43, 12).Value = "Segunda"
$ws.Cells.Item(43, 13).Value = 50
$ws.Cells.Item(43, 14).Value = 7000
$ws.Cells.Item(43, 15).Value = 7000
$ws.Cells.Item(43, 16).Value = 7000
$ws.Cells.Item(43, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(43, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(43, 19).Value = 389
$ws.Cells.Item(43, 20).Value = 18
